$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I ("I0") and J ("IF"), matching the formatting of
# the existing header cells (e.g. H1 - bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..37: I is always 1, J mirrors the existing H value for the row.
for ($r = 2; $r -le 37; $r++) {
    $hValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValue
}
